# BIOL8002 databases exercise: BRCA2 -> ACE2 rewrite, plus assorted wording
# fixes, as described by the commit's diff.
#
# Strategy: use Document.Content.Find.Execute with whole (unique) sentences
# as the search text and the corrected sentence as the replacement. This
# keeps each call unambiguous (no partial-word collisions) regardless of how
# Word happened to have split the original sentence across runs.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $find"
    }
}

# Title line
Replace-Text "Collate technical information about the BRCA2 gene in human genome" "Collate technical information about the ACE2 gene in human genome"

# NCBI search term
Replace-Text "“Homo sapiens” AND “BRCA2”" "“Homo sapiens” AND “ACE2”"

# Basic description question (reworded + gene swap)
Replace-Text "What is the basic description provided for the BRCA2 gene in terms of its function?" "What is the basic description provided for the ACE2 gene in humans?"

# Alternate names question
Replace-Text "What are the alternate names of the BRCA2 gene?" "What are the alternate names of the ACE2 gene?"

# Transcript/protein sequence question (reworded + gene swap)
Replace-Text "What is the transcript and protein sequences for the BRCA2 gene?" "What are the transcript and protein sequence IDs for the ACE2 gene?"

# Respective lengths question
Replace-Text "What are respective lengths of transcript and protein sequence for the BRCA2 gene?" "What are respective lengths of transcript and protein sequence for the ACE2 gene?"

# Genome Browser -> Genome Data Viewer
Replace-Text "Click on the Genome Browser to visualize the gene body." "Click on the Genome Data Viewer to visualize the gene body."

# Upstream/downstream genes question
Replace-Text "What are the genes immediately upstream and downstream of the BRCA2 gene?" "What are the genes immediately upstream and downstream of the ACE2 gene?"

# Location question (reworded: adds assembly version)
Replace-Text "What is the location of the BRCA2 gene in the human genome? " "What is the location of the ACE2 gene in the human genome assembly version GRCh38.p13? "

# Duplicate "chromosome" typo fix
Replace-Text "What chromosome the human X chromosome is homologous to in the following genomes" "What chromosome the human X is homologous to in the following genomes"

# Explore information / search bar (two BRCA2 occurrences in one sentence)
Replace-Text "Explore information about the BRCA2 gene by searching for BRCA2 in the search bar." "Explore information about the ACE2 gene by searching for ACE2 in the search bar."

# Select entry for human gene
Replace-Text "Select the entry for human BRCA2 gene." "Select the entry for human ACE2 gene."

# Alternate transcripts question
Replace-Text "How many alternate transcripts does BRCA2 gene have in the human genome?" "How many alternate transcripts does ACE2 gene have in the human genome?"

# Transcript types question
Replace-Text "What are different types of transcripts produced by the BRCA2 gene?" "What are different types of transcripts produced by the ACE2 gene?"

# Evolution question
Replace-Text "When did the BRCA2 gene evolve" "When did the ACE2 gene evolve"

# Copies of gene question
Replace-Text "How many copies of BRCA2 gene is found in " "How many copies of ACE2 gene is found in "

# Species with more than one copy question (also collapses an extra run/space)
Replace-Text "Are there any species with more than one copy of BRCA2 " "Are there any species with more than one copy of ACE2 "

# Ensembl Genes build number bump
Replace-Text "“Ensembl Genes 97” and “Human Genes”" "“Ensembl Genes 104” and “Human Genes”"

# Search for gene in human genome
Replace-Text "Search for BRCA2 gene in the human genome." "Search for ACE2 gene in the human genome."

# Compare the conservation of gene
Replace-Text "Compare the conservation of BRCA2 gene " "Compare the conservation of ACE2 gene "

# Bar plot expression patterns
Replace-Text "Click on the bar plot to see the expression patterns of the BRCA2 gene in various tissues in humans." "Click on the bar plot to see the expression patterns of the ACE2 gene in various tissues in humans."

# Most abundant expression question (fixes duplicated "the")
Replace-Text "What tissue has the the most abundant expression of BRCA2?" "What tissue has the most abundant expression of ACE2?"

# Compare expression bar plots sentence (reworded: generic "other genes from the pathway")
Replace-Text "Compare the expression bar plots of BRCA2 with N4BP2L1 (downstream of BRCA2) and FRY (upstream of BRCA2) in " "Compare the expression bar plots of ACE2 with other genes from the pathway in "

Write-Output "Done"
